$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - Ar condicionado (A/C)
$ws.Range("A1").Value = "ar"
$ws.Range("B1").Value = "A/C"
$ws.Range("C1").Value = 30
$ws.Range("D1").Value = $false
$ws.Range("F1").Value = $false

# Row 2 - Lampada
$ws.Range("A2").Value = "l2"
$ws.Range("B2").Value = "Lâmpada"
$ws.Range("C2").Value = 100
$ws.Range("D2").Value = $false
$ws.Range("F2").Value = $false

# Row 3 - Televisor
$ws.Range("A3").Value = "tv"
$ws.Range("B3").Value = "Televisor"
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = $false
$ws.Range("F3").Value = $false
